$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "added grades of lsn 9" - update attendance (col B) / assignment (col C)
# grades for several students (rows 16-27 of Sheet1)
$ws.Range("B16").Value = 26
$ws.Range("C16").Value = 35

$ws.Range("B18").Value = 35
$ws.Range("C18").Value = 35

$ws.Range("B19").Value = 30
$ws.Range("C19").Value = 15

$ws.Range("B20").Value = 27
$ws.Range("C20").Value = 25

$ws.Range("B23").Value = 35
$ws.Range("C23").Value = 15

$ws.Range("B26").Value = 28

# Update the visible view/selection to match where the new data was entered
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("C26").Select()
